$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 30: Question 1 under 06_CLCL module ---
$q30 = @'
SELECT CL.CLCL_ID, CL.CLCL_CUR_STS 
FROM CMC_CLCL_CLAIM CL,  CMC_SBSB_SUBSC SB 
WHERE CL.SBSB_CK = SB.SBSB_CK  
AND CL.CLCL_CUR_STS = '11' AND SB.SBSB_ID = '070700003' AND SB.SBSB_ORIG_EFF_DT<=GETDATE() AND SB.SBSB_MCTR_STS = 'ACTI';
'@
$ws.Range("C30").Value = $q30
$ws.Range("C30").WrapText = $true
$ws.Rows.Item(30).RowHeight = 52.5

# --- Row 31: Question 2 ---
$q31 = @'
SELECT CL.CLCL_ID, CL.CLCL_CL_TYPE, CL.CLCL_CL_SUB_TYPE
FROM CMC_CLCL_CLAIM CL, CMC_SBSB_SUBSC SB
WHERE CL.SBSB_CK = SB.SBSB_CK AND (CL.CLCL_CUR_STS = 11 OR CL.CLCL_CUR_STS = 15) 
AND SB.SBSB_ID = '070700003' AND SB.SBSB_ORIG_EFF_DT<=GETDATE()  AND  SB.SBSB_MCTR_STS='ACTI';
'@
$ws.Range("C31").Value = $q31
$ws.Range("C31").WrapText = $true
$ws.Rows.Item(31).RowHeight = 52.5

# --- Row 32: Question 3 ---
$q32 = @'
SELECT PR.PRPR_NAME, CD.IPCD_ID, CD.IDCD_ID 
FROM CMC_CLCL_CLAIM CL, CMC_CDML_CL_LINE CD, CMC_PRPR_PROV PR
WHERE CL.CLCL_ID = CD.CLCL_ID AND CD.PRPR_ID = PR.PRPR_ID AND CL.PRPR_ID = PR.PRPR_ID
AND CL.CLCL_ID = '072180000100';
'@
$ws.Range("C32").Value = $q32
$ws.Range("C32").WrapText = $true
$ws.Rows.Item(32).RowHeight = 52.5

# --- Row 33: Question 4 ---
$q33 = @'
SELECT pypy.PYPY_PAYOR_NAME, pyba.PYBA_BANK_NAME
FROM CMC_CLCK_CLM_CHECK clck, CMC_LOBD_LINE_BUS lobd, CMC_PYPY_PAYOR pypy, CMC_PYBA_BANK_ACCT pyba
WHERE clck.LOBD_ID = lobd.LOBD_ID
AND pypy.PYBA_ID = pyba.PYBA_ID
AND lobd.PYPY_ID = pypy.PYPY_ID 
AND clck.CLCL_ID = '072180000100' 
AND GETDATE() BETWEEN pypy.PYPY_EFF_DT AND pypy.PYPY_TERM_DT
'@
$ws.Range("C33").Value = $q33
$ws.Range("C33").WrapText = $true
$ws.Rows.Item(33).RowHeight = 91.5

# --- Row 34: Question 5 ---
$q34 = @'
SELECT CLCK.CKPY_REF_ID FROM CMC_CDML_CL_LINE CDML,CMC_CLCK_CLM_CHECK CLCK
WHERE CDML.CLCL_ID=CLCK.CLCL_ID
AND CDML.CLCL_ID='072180000100';
'@
$ws.Range("C34").Value = $q34
$ws.Range("C34").WrapText = $true
$ws.Rows.Item(34).RowHeight = 39.5

# --- Row 35: Question 6 (style/height unchanged - stays style s=6, default row height) ---
$q35 = "SELECT CLOV_AMT FROM CMC_CLOV_OVERPAY WHERE CLCL_ID = '072200000401';"
$ws.Range("C35").Value = $q35

# --- Row 36: Question 7 ---
$q36 = @'
SELECT ACPR.ACPR_RECOV_AMT, ACPR.ACPR_REF_ID
FROM CMC_ACPR_PYMT_RED ACPR, CMC_LOBD_LINE_BUS LOBD, CMC_CLCK_CLM_CHECK CLCK
WHERE ACPR.LOBD_ID = LOBD.LOBD_ID AND LOBD.LOBD_ID = CLCK.LOBD_ID AND CLCK.CLCL_ID = '072200000401'
'@
$ws.Range("C36").Value = $q36
$ws.Range("C36").WrapText = $true
$ws.Rows.Item(36).RowHeight = 39.5

# --- Row 37: Question 8 ---
$q37 = @'
SELECT CKPY.CKPY_REF_ID, ACPR.ACPR_RECOV_AMT
FROM CMC_CLCK_CLM_CHECK CLCK, CMC_LOBD_LINE_BUS LOBD, CMC_CKPY_PAYEE_SUM CKPY, CMC_ACRH_RED_HIST ACRH, CMC_ACPR_PYMT_RED ACPR
WHERE CLCK.LOBD_ID = LOBD.LOBD_ID AND LOBD.LOBD_ID = CKPY.LOBD_ID AND CKPY.CKPY_REF_ID = ACRH.CKPY_REF_ID AND ACRH.ACPR_REF_ID = ACPR.ACPR_REF_ID
AND CLCK.CLCL_ID = '072200000401'
'@
$ws.Range("C37").Value = $q37
$ws.Range("C37").WrapText = $true
$ws.Rows.Item(37).RowHeight = 52.5

# --- Update the view: scroll to/select B38, matching the saved cursor position ---
$ws.Range("B38").Select()
